# [ADD] TESTE COM LOGIN USANDO DB
#
# The original sheet listed two accounts:
#   Row 2: emervin / emersona7x@hotmail.com / Emerson Rafael / <hash> / ADM
#   Row 3: rbriggs / rbriggs@gmail.com      / Rebecca Briggs / <hash> / EDITORA
#
# This change drops the "emervin" (Emerson Rafael / ADM) test account,
# leaving only the "rbriggs" (Rebecca Briggs / EDITORA) row, which shifts
# up from row 3 to row 2. The mailto hyperlink that lived on B2 (the
# emervin row's email) goes away along with the row it was attached to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire second row (emervin / Emerson Rafael / ADM). Excel
# shifts row 3 (rbriggs / Rebecca Briggs / EDITORA) up into its place,
# and the sheet's dimension shrinks from A1:E3 to A1:E2 automatically.
$ws.Rows.Item(2).Delete()

# The hyperlink that used to sit on B2 belonged to the deleted account's
# email address; drop it so no stray <hyperlinks> entry (and its
# relationship) is left behind.
$ws.Range("B2").Hyperlinks.Delete()
